$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Row 87 - remove empty C87/D87 cells (they become absent, as in the target) ---
$ws.Range("C87").ClearContents()
$ws.Range("D87").ClearContents()

# --- Step 2: Append new rows 88-111 for case 21TRD09437 / Bunner ---
$newRows = @(
    @{Row=88; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F=""},
    @{Row=89; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F=""},
    @{Row=90; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F=""},
    @{Row=91; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Guilty"},
    @{Row=92; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Dismissed"},
    @{Row=93; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Guilty"},
    @{Row=94; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Guilty"},
    @{Row=95; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Dismissed"},
    @{Row=96; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Guilty"},
    @{Row=97; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F=""},
    @{Row=98; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F=""},
    @{Row=99; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F=""},
    @{Row=100; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Guilty"},
    @{Row=101; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Guilty"},
    @{Row=102; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Dismissed"},
    @{Row=103; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Dismissed"},
    @{Row=104; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Guilty"},
    @{Row=105; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Guilty"},
    @{Row=106; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Guilty"},
    @{Row=107; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Dismissed"},
    @{Row=108; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Guilty"},
    @{Row=109; A="21TRD09437"; B="Bunner"; C="DUS"; D="4510.11"; E="M1"; F="Guilty"},
    @{Row=110; A="21TRD09437"; B="Bunner"; C="1ST SPEED 1 YR SCHOOL >35MPHM4"; D="4511.21B1A"; E="M4"; F="Dismissed"},
    @{Row=111; A="21TRD09437"; B="Bunner"; C="RECKLESS OPERATION 1ST IN 1 YR"; D="4511.20"; E="MM"; F="Guilty"}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).NumberFormat = "@"
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    if ($r.F -ne "") {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
}
